# SQ1_Final_Results/sq1_metrics_rules.xlsx
# Replace the DeepSeek-R1 results (rows 2-11, 4 models x up to 3 groups)
# with the Gemma3 results (rows 2-7, 3 models x up to 2-3 groups).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing rows that belonged to deepseek_r1_14b / deepseek_r1_32b;
# only 6 data rows remain for the gemma3 models.
$ws.Rows("8:11").Delete()

# Columns: A=Model B=Group C=N D=V1_Tot E=V1_Act F=V2_Tot G=V2_Act H=V3_Tot
#          I=V3_Act J=Intv K=Intv_S L=Intv_P M=Intv_H N=FF O=Audit_Str
$rows = @(
    @{ Row=2; A="gemma3_4b";  B="Group_A"; C=998;  D=1;  E=1; F=36; G=36; H=0;   I=0; J=37;  K=0;  L=0; M=0; N=14.14; O="1|0" },
    @{ Row=3; A="gemma3_4b";  B="Group_B"; C=827;  D=42; E=0; F=0;  G=0;  H=99;  I=0; J=141; K=42; L=0; M=0; N=39.61; O="1|32" },
    @{ Row=4; A="gemma3_4b";  B="Group_C"; C=843;  D=53; E=0; F=0;  G=0;  H=115; I=0; J=168; K=53; L=0; M=0; N=37.42; O="0|43" },
    @{ Row=5; A="gemma3_12b"; B="Group_A"; C=999;  D=0;  E=0; F=89; G=89; H=0;   I=0; J=89;  K=0;  L=0; M=0; N=17.58; O="0|1" },
    @{ Row=6; A="gemma3_12b"; B="Group_B"; C=969;  D=10; E=0; F=22; G=2;  H=0;   I=0; J=32;  K=10; L=0; M=0; N=24.86; O="2|4" },
    @{ Row=7; A="gemma3_27b"; B="Group_A"; C=1000; D=0;  E=0; F=18; G=18; H=0;   I=0; J=18;  K=0;  L=0; M=0; N=32.67; O="0|0" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("N$n").Value = $r.N
    $ws.Range("O$n").Value = $r.O
}
